# Auto-generated Excel COM-interop script to apply the diff changes
# to the Cuchulainn Profits workbook (scheduled-runner update).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1250142.5
$ws.Range("I6").Value = 1666690
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 5000070
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -4999958
$ws.Range("N6").Value = -1724
$ws.Range("H17").Value = 2300
$ws.Range("I17").Value = 2300
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 6900
$ws.Range("L17").Value = $null
$ws.Range("N17").Value = 0
$ws.Range("M17").Value = -6732
$ws.Range("H20").Value = 2510.5
$ws.Range("I20").Value = 2510.5
$ws.Range("K20").Value = 2510.5
$ws.Range("M20").Value = -2280.5
$ws.Range("H35").Value = 2510.5
$ws.Range("I35").Value = 2510.5
$ws.Range("K35").Value = 2510.5
$ws.Range("M35").Value = -2131.5
$ws.Range("H55").Value = 660.1429000000001
$ws.Range("I55").Value = 606.75
$ws.Range("J55").Value = 731.3333
$ws.Range("K55").Value = 606.75
$ws.Range("L55").Value = 731.3333
$ws.Range("M55").Value = -392.75
$ws.Range("N55").Value = -1159.3333
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = $null
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = $null
$ws.Range("N77").Value = 0
$ws.Range("H98").Value = 7466.7144
$ws.Range("I98").Value = 1730.2
$ws.Range("K98").Value = 1730.2
$ws.Range("M98").Value = -232.2
$ws.Range("H122").Value = 7466.7144
$ws.Range("I122").Value = 1730.2
$ws.Range("K122").Value = 5190.6
$ws.Range("M122").Value = -2740.6
$ws.Range("H138").Value = 9465
$ws.Range("J138").Value = 14999.333
$ws.Range("L138").Value = 44997.999
$ws.Range("N138").Value = -55277.999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 16354
$ws.Range("J24").Value = 16354
$ws.Range("L24").Value = 16354
$ws.Range("N24").Value = -17102
$ws.Range("H32").Value = 3467.5
$ws.Range("I32").Value = 3161
$ws.Range("K32").Value = 3161
$ws.Range("M32").Value = -2874
$ws.Range("H53").Value = 6799
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 6799
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = $null
$ws.Range("M53").Value = 6799
$ws.Range("N53").Value = -8163
$ws.Range("H74").Value = 5250
$ws.Range("I74").Value = 6250
$ws.Range("J74").Value = 4250
$ws.Range("K74").Value = 6250
$ws.Range("L74").Value = 4250
$ws.Range("M74").Value = -5376
$ws.Range("N74").Value = -5998
$ws.Range("H77").Value = 5250
$ws.Range("I77").Value = 6250
$ws.Range("J77").Value = 4250
$ws.Range("K77").Value = 31250
$ws.Range("L77").Value = 21250
$ws.Range("M77").Value = -26882
$ws.Range("N77").Value = -29986
$ws.Range("H96").Value = 18998
$ws.Range("J96").Value = 18998
$ws.Range("L96").Value = 18998
$ws.Range("N96").Value = -24490
$ws.Range("H97").Value = 641.8570999999999
$ws.Range("I97").Value = 738.8
$ws.Range("J97").Value = 399.5
$ws.Range("K97").Value = 738.8
$ws.Range("L97").Value = 399.5
$ws.Range("M97").Value = -242.8
$ws.Range("N97").Value = -1391.5
$ws.Range("H100").Value = 16354
$ws.Range("J100").Value = 16354
$ws.Range("L100").Value = 16354
$ws.Range("N100").Value = -18518

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 30685.285
$ws.Range("J35").Value = 30685.285
$ws.Range("L35").Value = 30685.285
$ws.Range("N35").Value = -31305.285
$ws.Range("H58").Value = 46666
$ws.Range("J58").Value = 46666
$ws.Range("L58").Value = 46666
$ws.Range("N58").Value = -47254
$ws.Range("H60").Value = 100000
$ws.Range("J60").Value = 100000
$ws.Range("L60").Value = 100000
$ws.Range("N60").Value = -101198
$ws.Range("H94").Value = 1839
$ws.Range("I94").Value = 2748.5
$ws.Range("J94").Value = 20
$ws.Range("K94").Value = 2748.5
$ws.Range("L94").Value = 20
$ws.Range("M94").Value = -2297.5
$ws.Range("N94").Value = -922
$ws.Range("H105").Value = 4999
$ws.Range("I105").Value = 4999
$ws.Range("K105").Value = 4999
$ws.Range("M105").Value = -3252

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9000
$ws.Range("H34").Value = 9000
$ws.Range("H58").Value = 8033.2666
$ws.Range("I58").Value = 7599
$ws.Range("J58").Value = 8684.666999999999
$ws.Range("K58").Value = 7599
$ws.Range("L58").Value = 8684.666999999999
$ws.Range("M58").Value = -7396
$ws.Range("N58").Value = -9090.666999999999
$ws.Range("H62").Value = 2200
$ws.Range("J62").Value = 2200
$ws.Range("L62").Value = 2200
$ws.Range("N62").Value = -3448
$ws.Range("H65").Value = 2200
$ws.Range("J65").Value = 2200
$ws.Range("L65").Value = 11000
$ws.Range("N65").Value = -17240
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = $null
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = $null
$ws.Range("N77").Value = 0
$ws.Range("H88").Value = 20239.5
$ws.Range("J88").Value = 20239.5
$ws.Range("L88").Value = 20239.5
$ws.Range("N88").Value = -21051.5
$ws.Range("H91").Value = 20239.5
$ws.Range("J91").Value = 20239.5
$ws.Range("L91").Value = 20239.5
$ws.Range("N91").Value = -23047.5
$ws.Range("H96").Value = 16999.666
$ws.Range("J96").Value = 16999.666
$ws.Range("L96").Value = 16999.666
$ws.Range("N96").Value = -22491.666
$ws.Range("H136").Value = 8033.2666
$ws.Range("I136").Value = 7599
$ws.Range("J136").Value = 8684.666999999999
$ws.Range("K136").Value = 22797
$ws.Range("L136").Value = 26054.001
$ws.Range("M136").Value = -20247
$ws.Range("N136").Value = -31154.001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 575
$ws.Range("I11").Value = 600
$ws.Range("K11").Value = 1800
$ws.Range("M11").Value = -1660
$ws.Range("H23").Value = 309.5
$ws.Range("I23").Value = 500
$ws.Range("K23").Value = 1500
$ws.Range("M23").Value = -1265
$ws.Range("H29").Value = 300
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = $null
$ws.Range("M29").Value = 900
$ws.Range("N29").Value = -1454
$ws.Range("H46").Value = 100
$ws.Range("I46").Value = 100
$ws.Range("K46").Value = 300
$ws.Range("M46").Value = -209
$ws.Range("H131").Value = 3064.3333
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = $null
$ws.Range("N132").Value = 0

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 16448
$ws.Range("J95").Value = 16448
$ws.Range("L95").Value = 16448
$ws.Range("N95").Value = -21940
$ws.Range("H105").Value = 25990.666
$ws.Range("J105").Value = 25990.666
$ws.Range("L105").Value = 25990.666
$ws.Range("N105").Value = -32978.666
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530
$ws.Range("H132").Value = 6068.778
$ws.Range("I132").Value = 2723.8
$ws.Range("K132").Value = 8171.400000000001
$ws.Range("M132").Value = -5641.400000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 9499.75
$ws.Range("J104").Value = 9499.75
$ws.Range("L104").Value = 9499.75
$ws.Range("N104").Value = -16487.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = $null

